$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet to reflect the new date
$ws.Name = "Product-2023-09-12"

# 2. Fix header capitalization for C1 ("promo price" -> "Promo price")
$ws.Range("C1").Value = "Promo price"

# 3. Update date/time related text values for both data rows (row 2 and row 3)
$ws.Range("H2").Value = "12/09/2023 00:17"
$ws.Range("H3").Value = "12/09/2023 00:17"

# Keep these as plain text (ISO date strings look like real dates, and
# Excel would otherwise auto-convert them to date serials) before writing.
$ws.Range("I2:I3").NumberFormat = "@"
$ws.Range("I2").Value = "2023-09-12"
$ws.Range("I3").Value = "2023-09-12"

$ws.Range("J2").Value = "2023-09-12T00:17:58.4095525+01:00[Africa/Casablanca]"
$ws.Range("J3").Value = "2023-09-12T00:17:58.4095525+01:00[Africa/Casablanca]"

$ws.Range("L2").Value = "2023-09-12T00:17:58.4095525"
$ws.Range("L3").Value = "2023-09-12T00:17:58.4095525"

# 4. Column C widened slightly (auto-fit effect of the new, longer header
#    text "Promo price"). The host quantizes ColumnWidth to 1/6-character
#    pixel steps, so 11.8333 is the input that lands on the value closest
#    to the recorded target width (12.68359375 -> ~12.6667).
$ws.Columns.Item(3).ColumnWidth = 11.8333333
